$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.717.45'
$ws.Range("E2").Value = '  +0.74%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.319.72'
$ws.Range("E3").Value = '  +4.95%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.52'
$ws.Range("E5").Value = '  +3.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.39'
$ws.Range("E6").Value = '  +2.57%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.316.25'
$ws.Range("E8").Value = '  +5.10%  '
$ws.Range("E9").Value = '  +0.45%  '
$ws.Range("E10").Value = '  +2.48%  '
$ws.Range("E11").Value = '  +3.93%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.470'
$ws.Range("E12").Value = '  +2.17%  '
$ws.Range("E13").Value = '  +0.95%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.64'
$ws.Range("E14").Value = '  +1.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.861.09'
$ws.Range("E15").Value = '  +4.84%  '
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.314.04'
$ws.Range("E17").Value = '  +5.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.803.63'
$ws.Range("E18").Value = '  +0.88%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.86'
$ws.Range("E19").Value = '  +2.61%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '480.50'
$ws.Range("E20").Value = '  +0.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.18'
$ws.Range("E21").Value = '  +0.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.733'
$ws.Range("E22").Value = '  +4.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.11'
$ws.Range("E23").Value = '  +4.95%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.73'
$ws.Range("E24").Value = '  +5.48%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.68'
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("E26").Value = '  +0.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.77'
$ws.Range("E27").Value = '  +1.90%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.39'
$ws.Range("E28").Value = '  +2.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.14'
$ws.Range("E30").Value = '  +1.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.16'
$ws.Range("E31").Value = '  +1.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.83'
$ws.Range("E32").Value = '  +7.04%  '
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.54'
$ws.Range("E34").Value = '  +0.07%  '
$ws.Range("E35").Value = '  +2.98%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.02'
$ws.Range("E36").Value = '  +3.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.62'
$ws.Range("E37").Value = '  +0.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0745'
$ws.Range("E38").Value = '  +6.23%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0402'
$ws.Range("E39").Value = '  +3.20%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '432.75'
$ws.Range("E40").Value = '  +2.46%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.068.76'
$ws.Range("E41").Value = '  +4.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.75'
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.33'
$ws.Range("E43").Value = '  +0.50%  '
$ws.Range("E44").Value = '  +3.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.264'
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("E46").Value = '  +2.99%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.34'
$ws.Range("E47").Value = '  +3.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '36.18'
$ws.Range("E48").Value = '  +13.13%  '
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '125.52'
$ws.Range("E50").Value = '  +3.78%  '
$ws.Range("E51").Value = '  +0.73%  '
